$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35; this shifts the existing rows 35-54 down to 36-55
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new data record
$ws.Range("A35").Value = 3
$ws.Range("B35").Value = "Femacal de La Calera"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 44755
$ws.Range("E35").Value = 5
$ws.Range("F35").Value = 100112035
$ws.Range("G35").Value = "Bruselas (repollito)"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 100
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 16000
$ws.Range("M35").Value = 15550
$ws.Range("N35").Value = "$/malla 15 kilos"
$ws.Range("O35").Value = "Provincia de Quillota"
$ws.Range("P35").Value = 1037
$ws.Range("Q35").Value = 15
$ws.Range("R35").Value = "Hortaliza"
